$wb = $excel.ActiveWorkbook

# --- Sheet 1 ("土地"): a new land-property record was found, add it as
#     the new second data row (index 14), pushing the existing rows down.
$ws1 = $wb.Worksheets.Item(1)

# Insert a blank row above the current row 2, shifting rows 2-6 down to 3-7,
# then seed it from the row that is now directly below (keeps the A-column
# index style consistent with the rest of the table).
$ws1.Rows.Item(2).Insert()
$ws1.Range("A3:Q3").Copy($ws1.Range("A2:Q2"))

$ws1.Range("A2").Value = 14
$ws1.Range("B2").Value = "南投縣仁愛郷廬山段06760000地號"
$ws1.Range("C2").Value = 17197
$ws1.Range("D2").Value = "2分之1"
$ws1.Range("E2").Value = "孔文吉"
$ws1.Range("F2").NumberFormat = "@"
$ws1.Range("F2").Value = "95年10月03日"
$ws1.Range("G2").Value = "受贈"
$ws1.Range("H2").ClearContents()
$ws1.Range("I2").Value = "land"
$ws1.Range("J2").Value = "normal"
$ws1.Range("K2").NumberFormat = "@"
$ws1.Range("K2").Value = "2012-04-26"
$ws1.Range("L2").Value = "孔文吉"
$ws1.Range("M2").Value = 1312
$ws1.Range("N2").Value = "tmpfed71"
$ws1.Range("O2").Value = 14
$ws1.Range("P2").Value = 0.5
$ws1.Range("Q2").Value = 8598.5

# --- Sheet 2 ("汽車"): the single car record now also gets an explicit
#     index row (34), duplicating the same data onto row 2.
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A2").Value = 34
$ws2.Range("B2").Value = "中華自小貨"
$ws2.Range("C2").Value = 1094
$ws2.Range("D2").Value = "陳秋月"
$ws2.Range("E2").NumberFormat = "@"
$ws2.Range("E2").Value = "98年01月31闩"
$ws2.Range("F2").Value = "受贈"
$ws2.Range("G2").Value = 50000
